# Append a new paragraph (with a grammar proof-error marker around the
# word "report") plus one trailing empty paragraph to the end of the
# document body, right before the section properties — matching:
#
#   <w:p>
#     <w:r><w:t xml:space="preserve">This is a document that is included in the final security </w:t></w:r>
#     <w:proofErr w:type="gramStart"/>
#     <w:r><w:t>report</w:t></w:r>
#     <w:proofErr w:type="gramEnd"/>
#   </w:p>
#   <w:p/>

$d = $word.ActiveDocument

# Collapsed range sitting at the very end of the document's main story.
$endPos = $d.Content.End
$insertionRange = $d.Range($endPos, $endPos)

# A self-contained WordProcessingML package fragment: InsertXML replaces
# the (collapsed, i.e. zero-length) target range with this content,
# giving us exact control over the run/proofErr structure -- including
# emitting a truly empty trailing <w:p/> (no run at all).
$xmlFragment = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:t xml:space="preserve">This is a document that is included in the final security </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>report</w:t></w:r><w:proofErr w:type="gramEnd"/></w:p><w:p/></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

$insertionRange.InsertXML($xmlFragment)
